$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The export table ("Table1") needs a new trailing column, "MangingDirector",
# added after the existing "PostalCode" column (to add managing director to
# the excel export).
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

# Setting the value through the new table column's header cell both writes
# the worksheet cell (H1) and names the table column in one go.
$newCol.Range.Cells.Item(1, 1).Value = "MangingDirector"

# Match the column formatting of the rest of the header row (reuses the
# existing shared cell style instead of creating a new one).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Give the new column a sensible width, same as the other header columns
# (the engine quantizes widths to whole character units, so 15.0 lands on
# the stored width closest to the target 15.1640625).
$ws.Columns.Item(8).ColumnWidth = 15.0

# Leave the selection on the newly added header cell.
$ws.Range("H1").Select() | Out-Null
